$d = $word.ActiveDocument

# --- Locate the relevant paragraphs robustly (by content) instead of hard-coded indices ---
$fr36Index = -1
$outputIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($fr36Index -eq -1 -and $t -like "FR3.6*") {
        $fr36Index = $i
    }
    if ($fr36Index -ne -1 -and $i -gt $fr36Index -and $t -like "*输出：*") {
        $outputIndex = $i
        break
    }
}

if ($fr36Index -eq -1 -or $outputIndex -eq -1) {
    Write-Output "ERROR: could not locate paragraphs (fr36=$fr36Index output=$outputIndex)"
} else {
    Write-Output "fr36Index=$fr36Index outputIndex=$outputIndex"

    # --- 1. Remove the _GoBack bookmark from its current location (end of the FR3.6 title paragraph) ---
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    # --- 2. Re-add the _GoBack bookmark, collapsed, right after the text of the "输出" paragraph ---
    #     (immediately before its paragraph mark). A collapsed range placed exactly one character
    #     before a paragraph mark is mishandled by Bookmarks.Add in this runtime, so we work around
    #     it by temporarily appending a sentinel character, anchoring the bookmark next to it, then
    #     removing the sentinel again (the bookmark stays correctly anchored to the real text).
    $outPara = $d.Paragraphs.Item($outputIndex)
    $tail = $outPara.Range.Duplicate
    $tail.Collapse(0)
    $tail.MoveEnd(1, -1)
    $tail.InsertAfter("@@@SENTINEL@@@")

    $outPara2 = $d.Paragraphs.Item($outputIndex)
    $bmPos = $outPara2.Range.Duplicate
    $bmPos.Collapse(0)
    $bmPos.MoveEnd(1, -1)
    $bmPos.MoveEnd(1, -("@@@SENTINEL@@@".Length))
    $bmPos.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $bmPos)

    # remove the sentinel text again
    $sentinelRange = $d.Content.Find
    $found = $d.Range($d.Paragraphs.Item($outputIndex).Range.Start, $d.Paragraphs.Item($outputIndex).Range.End)
    $found.Find.Execute("@@@SENTINEL@@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $found.Delete()

    # --- 3. Strip the paragraph-mark run formatting (<w:rPr> inside <w:pPr>) on the "输出" paragraph ---
    $outPara3 = $d.Paragraphs.Item($outputIndex)
    $outPara3.Range.ParagraphFormat.Reset()
    $outPara3.Range.ParagraphFormat.FirstLineIndent = 21
